$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle "London Borough" -> "Recommendation for London Tourist" ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Recommendation for London Tourist"

# --- Slide 2: Content placeholder, first paragraph gets an added sentence ---
$s2 = $p.Slides.Item(2)
$tf = $s2.Shapes.Item(2).TextFrame
$para1 = $tf.TextRange.Paragraphs(1, 1)

# Extend the existing run's text in place (keeps same run/formatting) so it
# reads "...boroughs of London. So that recommendations can be made to a tourist "
$run1Range = $para1.Characters(1, 114)
$run1Range.Text = "The intention of this project is to perform a neighborhood analysis on the venues of all the boroughs of London. So that recommendations can be made to a tourist "

# Append the remaining words as a new trailing run: "visiting London."
$para1Again = $tf.TextRange.Paragraphs(1, 1)
$endPoint = $para1Again.Characters($para1Again.Length, 0)
$null = $endPoint.InsertAfter("visiting London.")
